$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('STF_ENG, Abelanes Maria Teresa Robles', 'Module Leader'),
    @('STF_ENG, Azhar Bin Moh', 'LMS TSO'),
    @('STF_ENG, Chow Kok Peng', 'Module Leader'),
    @('STF_ENG, Ho Kuan Tat Jackson', 'Module Leader'),
    @('STF_ENG, Low Poh Lian', 'Lecturer (Teaching Team)'),
    @('STF_ENG, Lye Sau Lin', 'Module Leader'),
    @('STF_ENG, Tan Geok Ling', 'Module Leader'),
    @('STF_ENG, Wee May Lin', 'Module Leader'),
    @('STF_ENG, Yan Seow Chiang', 'Module Leader'),
    @('STF_ENG, Yeo Sock Kiang May', 'LMS TSO'),
    @('STU_ENG, ANG JIA JUN', 'Learner'),
    @('STU_ENG, CHIM JUN XI LEROY GABRIEL', 'Learner'),
    @('STU_ENG, DHANI ARJUNA PUTRA BIN SHAZALI', 'Learner'),
    @('STU_ENG, HUAN JIAX CHOU', 'Learner'),
    @('STU_ENG, JASON PHUA KAI XIANG', 'Learner'),
    @('STU_ENG, JOHNSON NG CHONG SHENG', 'Learner'),
    @('STU_ENG, KALLANAYAKANAHALLI KANISHK PRABHAKAR', 'Learner'),
    @('STU_ENG, LI HAOYANG', 'Learner'),
    @('STU_ENG, MAHAESWAR S/O RAJENDRA PRASAD', 'Learner'),
    @('STU_ENG, MOWER PRISCILLA LE XUAN', 'Learner'),
    @('STU_ENG, MUHAMMAD ILHAM GONZAGA', 'Learner'),
    @('STU_ENG, PANG ZI JIAN', 'Learner'),
    @('STU_ENG, SOH WEI QUAN ROYDON', 'Learner'),
    @('STU_ENG, TAN KE MING', 'Learner'),
    @('STU_ENG, TEO YONG JIE ASHER', 'Learner'),
    @('STU_ENG, VINAYA SENTHILKUMAR', 'Learner'),
    @('STU_ENG, WANG PENGYU', 'Learner'),
    @('STU_ENG, WANG RUI ZHE EDDRICK', 'Learner'),
    @('STU_ENG, YONG KAH HAO MARCUS', 'Learner'),
    @('STU_ENG, YUSUF RAMI NAZZAL', 'Learner'),
    @('STU_ENG, ZAW YE YINT', 'Learner'),
    @('STU_ENG, ZHANG RAN', 'Learner')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
}
